$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-05 20:47:48"
$ws.Range("O2").Value = "-1.0 °C"
$ws.Range("E3").Value = "2026-02-05 20:47:50"
$ws.Range("E4").Value = "2026-02-05 20:47:53"
$ws.Range("O4").Value = "11.6 °C"
$ws.Range("E5").Value = "2026-02-05 20:47:56"
$ws.Range("J5").Value = "989.8 hPa"
$ws.Range("E6").Value = "2026-02-05 20:47:58"
$ws.Range("L6").Value = "36.7 km/h - 295º 20:22 TU"
$ws.Range("O6").Value = "13.0 °C"
$ws.Range("E7").Value = "2026-02-05 20:48:01"
$ws.Range("H7").Value = "'79%"
$ws.Range("J7").Value = "991.5 hPa"
$ws.Range("E8").Value = "2026-02-05 20:48:03"
$ws.Range("O8").Value = "9.0 °C"
$ws.Range("E9").Value = "2026-02-05 20:48:06"
$ws.Range("M9").Value = "6.9 °C 20:20 TU"
$ws.Range("E10").Value = "2026-02-05 20:48:08"
$ws.Range("H10").Value = "'90%"
$ws.Range("O10").Value = "7.9 °C"
$ws.Range("E11").Value = "2026-02-05 20:48:10"
$ws.Range("J11").Value = "994.6 hPa"
$ws.Range("O11").Value = "0.8 °C"
$ws.Range("E12").Value = "2026-02-05 20:48:13"
$ws.Range("H12").Value = "'86%"
$ws.Range("O12").Value = "10.3 °C"
$ws.Range("E13").Value = "2026-02-05 20:48:16"
$ws.Range("E14").Value = "2026-02-05 20:48:18"
$ws.Range("I14").Value = "7.6 mm"
$ws.Range("E15").Value = "2026-02-05 20:48:20"
$ws.Range("H15").Value = "'80%"
$ws.Range("J15").Value = "990.3 hPa"
$ws.Range("O15").Value = "8.6 °C"
$ws.Range("E16").Value = "2026-02-05 20:48:23"
$ws.Range("O16").Value = "3.9 °C"
$ws.Range("E17").Value = "2026-02-05 20:48:25"
$ws.Range("J17").Value = "995.0 hPa"
$ws.Range("K17").Value = "1.9 MJ/m2"
$ws.Range("O17").Value = "1.0 °C"
$ws.Range("E18").Value = "2026-02-05 20:48:28"
$ws.Range("I18").Value = "2.4 mm"
$ws.Range("E19").Value = "2026-02-05 20:48:31"
$ws.Range("E20").Value = "2026-02-05 20:48:33"
$ws.Range("E21").Value = "2026-02-05 20:48:36"
$ws.Range("O21").Value = "6.5 °C"
$ws.Range("E22").Value = "2026-02-05 20:48:38"
$ws.Range("H22").Value = "'87%"
$ws.Range("O22").Value = "9.0 °C"
$ws.Range("E23").Value = "2026-02-05 20:48:41"
$ws.Range("E24").Value = "2026-02-05 20:48:43"
$ws.Range("O24").Value = "10.5 °C"
$ws.Range("E25").Value = "2026-02-05 20:48:46"
$ws.Range("E26").Value = "2026-02-05 20:48:48"
$ws.Range("E27").Value = "2026-02-05 20:48:51"
$ws.Range("E28").Value = "2026-02-05 20:48:53"
$ws.Range("H28").Value = "'94%"
$ws.Range("J28").Value = "992.8 hPa"
$ws.Range("M28").Value = "7.3 °C 20:25 TU"
$ws.Range("O28").Value = "2.7 °C"
$ws.Range("E29").Value = "2026-02-05 20:48:56"
$ws.Range("H29").Value = "'79%"
$ws.Range("O29").Value = "9.4 °C"
$ws.Range("E30").Value = "2026-02-05 20:48:58"
$ws.Range("E31").Value = "2026-02-05 20:49:01"
$ws.Range("J31").Value = "994.3 hPa"
$ws.Range("E32").Value = "2026-02-05 20:49:03"
$ws.Range("E33").Value = "2026-02-05 20:49:06"
$ws.Range("H33").Value = "'85%"
$ws.Range("O33").Value = "9.5 °C"
$ws.Range("E34").Value = "2026-02-05 20:49:08"
$ws.Range("H34").Value = "'95%"
$ws.Range("O34").Value = "4.2 °C"
$ws.Range("E35").Value = "2026-02-05 20:49:10"
$ws.Range("G35").Value = "200 cm"
$ws.Range("I35").Value = "5.4 mm"
$ws.Range("E36").Value = "2026-02-05 20:49:13"
